$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 437, shifting existing rows 437:475 down to 438:476.
$ws.Rows.Item(437).Insert()

# Populate the newly-inserted row 437 with the new data record.
$ws.Range("A437").Value = 3
$ws.Range("B437").Value = "Femacal de La Calera"
$ws.Range("C437").Value = "Coquimbo"
$ws.Range("D437").Value = 45106
$ws.Range("E437").Value = 5
$ws.Range("F437").Value = 100112001
$ws.Range("G437").Value = "Berenjena"
$ws.Range("H437").Value = "Sin especificar"
$ws.Range("I437").Value = "Primera"
$ws.Range("J437").Value = 130
$ws.Range("K437").Value = 5500
$ws.Range("L437").Value = 6000
$ws.Range("M437").Value = 5692
$ws.Range("N437").Value = "$/caja 60 unidades"
$ws.Range("O437").Value = "Región de Arica y Parinacota"
$ws.Range("P437").Value = 95
$ws.Range("Q437").Value = 60
$ws.Range("R437").Value = "Hortaliza"
